$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E holds plain-text dates like "2024-10-08" which Excel would
# otherwise auto-convert to a date serial on assignment. Pre-format the
# target cells as Text so the literal string is preserved, matching the
# existing rows (which are stored as literal text, not dates).
$ws.Range("E339:E342").NumberFormat = "@"

# Row 339
$ws.Range("A339").Value = "2024-10-08 21:15:41"
$ws.Range("B339").Value = "check_availability"
$ws.Range("C339").Value = "https://example.com/reservation"
$ws.Range("D339").Value = 'Checked availability: (''Checked availability: Availability confirmed'', ''Data saved to Excel file at ExportedFiles\excelFiles\check_availability.xlsx.'', ''HTML file saved and updated at ExportedFiles\htmlFiles\check_availability.html.'')'
$ws.Range("E339").Value = "2024-10-08"
$ws.Range("F339").Value = "21:15:41"

# Row 340
$ws.Range("A340").Value = "2024-10-08 21:18:32"
$ws.Range("B340").Value = "check_availability"
$ws.Range("C340").Value = "https://example.com/reservation"
$ws.Range("D340").Value = 'Checked availability: (''Checked availability: Availability confirmed'', ''Data saved to Excel file at ExportedFiles\excelFiles\check_availability.xlsx.'', ''HTML file saved and updated at ExportedFiles\htmlFiles\check_availability.html.'')'
$ws.Range("E340").Value = "2024-10-08"
$ws.Range("F340").Value = "21:18:32"

# Row 341
$ws.Range("A341").Value = "2024-10-08 21:21:53"
$ws.Range("B341").Value = "check_availability"
$ws.Range("C341").Value = "https://example.com/reservation"
$ws.Range("D341").Value = 'Checked availability: (''Checked availability: Availability confirmed'', ''Data saved to Excel file at ExportedFiles\excelFiles\check_availability.xlsx.'', ''HTML file saved and updated at ExportedFiles\htmlFiles\check_availability.html.'')'
$ws.Range("E341").Value = "2024-10-08"
$ws.Range("F341").Value = "21:21:53"

# Row 342
$ws.Range("A342").Value = "2024-10-08 21:24:54"
$ws.Range("B342").Value = "check_availability"
$ws.Range("C342").Value = "https://example.com/reservation"
$ws.Range("D342").Value = 'Checked availability: (''Checked availability: Availability confirmed'', ''Data saved to Excel file at ExportedFiles\excelFiles\check_availability.xlsx.'', ''HTML file saved and updated at ExportedFiles\htmlFiles\check_availability.html.'')'
$ws.Range("E342").Value = "2024-10-08"
$ws.Range("F342").Value = "21:24:54"
